$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.698.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.548.25'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '198.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '588.19'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.615'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.630'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.23'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000290'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.36'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '683.24'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +15.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.104.33'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.677.27'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.545.59'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.48'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.97%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.66'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.973'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.87'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '108.27'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.93%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.42'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.96'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.73'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.70'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.39'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.96'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '62.37'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.808.90'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0819'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.70'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.97'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '501.10'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.92%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.374'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '34.99'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0461'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.96'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.60%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.138'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.41'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +21.08%  '
